# Auto update Excel log
# Appends the latest Proximity (ENTER/EXIT) and Camera (Image Captured) sensor
# events to the log sheets, mirroring the SeniorConnect sensor pipeline.

$wb = $excel.ActiveWorkbook

# --- Proximity sheet: Living Room Main Door ENTER/EXIT pair -----------------
$proximity = $wb.Worksheets.Item("Proximity")
$pRow = $proximity.UsedRange.Rows.Count + 1

$proximity.Cells.Item($pRow, 1).Value = "'2026-02-01"
$proximity.Cells.Item($pRow, 1).Style = "Normal"
$proximity.Cells.Item($pRow, 2).Value = "13:38:23"
$proximity.Cells.Item($pRow, 3).Value = "13:00"
$proximity.Cells.Item($pRow, 4).Value = "Living Room Main Door"
$proximity.Cells.Item($pRow, 5).Value = "ENTER"
$proximity.Cells.Item($pRow, 6).Value = "User ENTERED Living Room Main Door"
$pRow++

$proximity.Cells.Item($pRow, 1).Value = "'2026-02-01"
$proximity.Cells.Item($pRow, 1).Style = "Normal"
$proximity.Cells.Item($pRow, 2).Value = "13:38:48"
$proximity.Cells.Item($pRow, 3).Value = "13:00"
$proximity.Cells.Item($pRow, 4).Value = "Living Room Main Door"
$proximity.Cells.Item($pRow, 5).Value = "EXIT"
$proximity.Cells.Item($pRow, 6).Value = "User EXITED Living Room Main Door"

# --- Camera sheet: matching "Image Captured" events --------------------------
$camera = $wb.Worksheets.Item("Camera")
$cRow = $camera.UsedRange.Rows.Count + 1

$camera.Cells.Item($cRow, 1).Value = "'2026-02-01"
$camera.Cells.Item($cRow, 1).Style = "Normal"
$camera.Cells.Item($cRow, 2).Value = "13:38:24"
$camera.Cells.Item($cRow, 3).Value = "13:00"
$camera.Cells.Item($cRow, 4).Value = "Living Room Main Door"
$camera.Cells.Item($cRow, 5).Value = "Image Captured"
$camera.Cells.Item($cRow, 6).Value = "Active"
$cRow++

$camera.Cells.Item($cRow, 1).Value = "'2026-02-01"
$camera.Cells.Item($cRow, 1).Style = "Normal"
$camera.Cells.Item($cRow, 2).Value = "13:38:48"
$camera.Cells.Item($cRow, 3).Value = "13:00"
$camera.Cells.Item($cRow, 4).Value = "Living Room Main Door"
$camera.Cells.Item($cRow, 5).Value = "Image Captured"
$camera.Cells.Item($cRow, 6).Value = "Active"
